# Rename the district "Lake Minara" -> "Lago Minara" throughout column G
# (Distrito de residencia), preserving each cell's original letter case
# ("Lake Minara" -> "Lago Minara", "LAKE MINARA" -> "LAGO MINARA").
#
# A case-SENSITIVE whole-column Find/Replace is used (MatchCase:=$true) so
# the differently-cased "LAKE MINARA" / "Lake Minara" variants each map to
# their corresponding "LAGO MINARA" / "Lago Minara" replacement, and other
# district values that happen to share letters (e.g. "L Minara", "L MINARA",
# "Lakeside") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colG = $ws.Columns.Item(7)

# xlWhole = 1 (match entire cell contents), LookAt:=xlWhole, MatchCase:=$true
[void]$colG.Replace("LAKE MINARA", "LAGO MINARA", -4163, 1, $true)
[void]$colG.Replace("Lake Minara", "Lago Minara", -4163, 1, $true)

Write-Output "Replaced 'Lake Minara' / 'LAKE MINARA' with 'Lago Minara' / 'LAGO MINARA' in column G"
